$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "103"
$ws.Range("B2").Value = "akshay"
$ws.Range("C2").Value = "2026-01-29"
$ws.Range("D2").Value = "18:50:51"
$ws.Range("E2").Value = "Present"
